$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.624.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.290.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.86%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'96.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'268.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.55%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -1.92%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.83%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'45.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.86%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.67%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.75%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.635.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.89%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.00%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.294.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.95%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'43.609.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.88%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.83%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'72.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.39%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +12.36%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'232.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.70%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -6.49%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.07%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.92%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'11.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.70%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +3.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'40.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.72%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.93%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'175.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.31%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'22.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.92%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0900"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.46%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.24%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.68%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.71%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.39%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +2.58%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.34%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'MultiversX"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'66.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +7.81%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Celestia"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'12.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'ARBITRUM"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.64%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'8.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.00%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -4.31%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.56%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'97.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.89%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.43%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +9.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.185"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +8.46%  "
$ws.Range("E51").Style = "Normal"
